$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "57.655.19"
$ws.Range("E2").Value = "  +0.89%  "

$ws.Range("D3").Value = "2.361.57"
$ws.Range("E3").Value = "  +1.71%  "

$ws.Range("E4").Value = "  +0.04%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "522.21"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +0.32%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "136.76"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +1.47%  "

$ws.Range("E7").Value = "  -0.10%  "

$ws.Range("E8").Value = "  +0.20%  "

$ws.Range("E9").Value = "  -0.76%  "

$ws.Range("E10").Value = "  +5.31%  "

$ws.Range("E11").Value = "  -0.95%  "

$ws.Range("E12").Value = "  -0.17%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "24.39"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +0.95%  "

$ws.Range("D14").Value = "2.781.91"
$ws.Range("E14").Value = "  +1.76%  "

$ws.Range("D15").Value = "57.667.87"
$ws.Range("E15").Value = "  +1.40%  "

$ws.Range("E16").Value = "  +0.01%  "

$ws.Range("D17").Value = "2.366.71"
$ws.Range("E17").Value = "  +1.43%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "10.66"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +0.59%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "331.36"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +2.89%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.25"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -1.20%  "

$ws.Range("E21").Value = "  +1.02%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.998"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -0.13%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "61.42"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +0.17%  "

$ws.Range("E24").Value = "  +3.84%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.994"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +0.07%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "8.39"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +8.78%  "

$ws.Range("E27").Value = "  +8.32%  "

$ws.Range("D28").Value = "0.0₃0750"
$ws.Range("E28").Value = "  +1.38%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "170.28"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -0.86%  "

$ws.Range("E30").Value = "  +0.44%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.30"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +0.08%  "

$ws.Range("E32").Value = "  +0.78%  "

$ws.Range("E33").Value = "  +0.02%  "

$ws.Range("E34").Value = "  +2.81%  "

$ws.Range("E35").Value = "  -0.26%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.930"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -2.33%  "

$ws.Range("E37").Value = "  -0.31%  "

$ws.Range("E38").Value = "  +3.96%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "38.65"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +2.94%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "151.38"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +7.38%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.385"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +0.51%  "

$ws.Range("E42").Value = "  +1.68%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.34"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +2.61%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "283.27"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +1.08%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0942"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +1.17%  "

$ws.Range("E46").Value = "  -0.40%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.566"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +1.38%  "

$ws.Range("B48").Value = "VeChain"
$ws.Range("C48").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0222"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +2.21%  "

$ws.Range("B49").Value = "InjectiveProtocol"
$ws.Range("C49").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "18.32"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +5.34%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "17.70"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +4.07%  "

$ws.Range("E51").Value = "  -0.69%  "
